$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking text stored as inlineStr/text cells.
# Use an apostrophe prefix to force text entry, then reset the style so no
# extra number-format / quote-prefix style sticks to the cell (matches source,
# which has no 's' attribute on these data cells).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "246.84"
Set-TextValue "D3" "21.75"
Set-TextValue "D4" "5.455"
Set-TextValue "D5" "0.05697"
Set-TextValue "D6" "3.371"
Set-TextValue "D7" "0.8006"
Set-TextValue "D8" "1.034"
Set-TextValue "D9" "0.1431"
Set-TextValue "D10" "0.07298"
Set-TextValue "D11" "0.03164"
Set-TextValue "D12" "0.02935"
Set-TextValue "D13" "0.09287"
Set-TextValue "D14" "0.001642"
Set-TextValue "D15" "3.207"
Set-TextValue "D16" "0.04708"
Set-TextValue "D17" "0.0005899"
$ws.Range("E17").Value = "16OneONE"
Set-TextValue "D18" "0.006399"
Set-TextValue "D19" "0.005045"
$ws.Range("E19").Value = "18HotbitTokenHTB"
Set-TextValue "D20" "0.001046"
Set-TextValue "D21" "0.0001500"
Set-TextValue "D22" "0.0003199"
Set-TextValue "D23" "3.801"
Set-TextValue "D24" "6.430"
Set-TextValue "D25" "2.088"
Set-TextValue "D26" "0.3288"
Set-TextValue "D40" "0.04084"
Set-TextValue "D41" "0.006919"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
Set-TextValue "D42" "0.003499"
Set-TextValue "D44" "0.008086"
Set-TextValue "D45" "0.00005835"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.6823"
Set-TextValue "D48" "0.01010"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.01010"
